# daily auto push: 2026-01-28 13:54 UTC
# Insert a new data row for 2026/01/28 03:00-ish slot ("19") into the
# daily ranking log sheet. The new row becomes row 739 and every row that
# used to be 739..780 shifts down by one (740..781).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 739..780 down by one to make room for the new row.
$ws.Rows.Item(739).Insert()

# Fill in the newly inserted row.
# Column A holds a date formatted as plain text (e.g. "2026/01/28"), not a
# real date value, so we enter it with a leading apostrophe (forces text)
# and then strip the resulting "quote prefix" formatting so the cell ends
# up as a plain, unstyled text cell exactly like the rest of the sheet.
$ws.Range("A739").Formula = "'2026/01/28"
$ws.Range("A739").ClearFormats()

$ws.Range("B739").Value = "水"
$ws.Range("C739").Value = 19
$ws.Range("D739").Value = 201
